# Apply "last changes to v1.8.2" updates to the workbook.
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# 1. Bump version number
$wsMeta.Range("B3").Value = "1.8.2"

# 2. Update the publication date/time
$wsMeta.Range("B8").Value = "2023-09-01T14:45:29-04:00"

# 3. Populate the previously-empty invariants cell for the root Extension row
#    with the same constraint text already present on Extension.extension (AJ3).
$invariantText = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
$wsElem.Range("AJ1").Value = $invariantText
